## Adds two new worksheets ("view_departments" and "edit_departments") to the
## workbook, mirroring the existing "view_centers" / "edit_centers" sheets,
## updates the department data on the new "edit_departments" sheet, and
## leaves "edit_departments" as the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- Create "view_departments" as a copy of "view_centers" ---------------
$wsViewCenters = $wb.Worksheets.Item("view_centers")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsViewCenters.Copy($null, $lastSheet)
$wsViewDept = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsViewDept.Name = "view_departments"

# --- Create "edit_departments" as a copy of "edit_centers" ---------------
$wsEditCenters = $wb.Worksheets.Item("edit_centers")
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEditCenters.Copy($null, $lastSheet2)
$wsEditDept = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsEditDept.Name = "edit_departments"

# --- Update the department codes / names on "edit_departments" -----------
$wsEditDept.Range("B2").Value = "DEPT_updated_0081"
$wsEditDept.Range("B2").Style = "Normal"
$wsEditDept.Range("B3").Value = "DEPT_updated_0082"
$wsEditDept.Range("B3").Style = "Normal"
$wsEditDept.Range("B4").Value = "DEPT_updated_0083"
$wsEditDept.Range("B4").Style = "Normal"
$wsEditDept.Range("B5").Value = "DEPT_updated_0084"
$wsEditDept.Range("B5").Style = "Normal"

$wsEditDept.Range("C2").Value = "Dept_name_udpated_99"
$wsEditDept.Range("C3").Value = "Dept_name_udpated_100"
$wsEditDept.Range("C4").Value = "Dept_name_udpated_101"
$wsEditDept.Range("C5").Value = "Dept_name_udpated_102"

# --- Restore per-sheet selections, matching the committed workbook -------
[void]$wsViewDept.Range("E13").Select()
[void]$wsEditDept.Range("H9").Select()

# --- Deselect the "add_new_locations" sheet as the "tabSelected" sheet ---
# (already handled: activating the new sheets above moves tabSelected away
# from add_new_locations to edit_departments automatically)
